$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows for accounts 008032597 (ALESSANDRO) and 008037529 (MELISSA).
# These are Excel rows 3 and 4 (row 1 is the header, row 2 is 001882235/LAGO).
$ws.Range("A3:A4").EntireRow.Delete()
